$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old Row 7 values
$ws.Range("D2").Value = 44559
$ws.Range("K2").Value = 'Modesto'
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("S2").Value = 1083
$ws.Range("T2").Value = 18

# Row 3 <- old Row 8 values
$ws.Range("D3").Value = 44559
$ws.Range("K3").Value = 'Modesto'
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 18

# Row 7 <- old Row 9 values
$ws.Range("D7").Value = 44187
$ws.Range("K7").Value = 'Dina'
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 15500
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 861

# Row 8 <- old Row 2 values
$ws.Range("D8").Value = 44545
$ws.Range("K8").Value = 'Castle Brite'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 100
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 18500
$ws.Range("Q8").Value = '$/caja 15 kilos'
$ws.Range("S8").Value = 1233
$ws.Range("T8").Value = 15

# Row 9 <- old Row 3 values
$ws.Range("D9").Value = 44545
$ws.Range("K9").Value = 'Castle Brite'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 17000
$ws.Range("Q9").Value = '$/caja 15 kilos'
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1133
$ws.Range("T9").Value = 15
